# Atualizacao de bases das ligas (Mexico Liga de Expansion)
# - updates existing rows 217-219 with corrected match data
# - appends new rows 220-223 (new matches), extending dimension to A1:AC223
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 217-219 ---
# Row 217
$ws.Cells.Item(217, 1).Value = 215
$ws.Cells.Item(217, 2).Value = 7640652
$ws.Cells.Item(217, 3).Value = "Mexico Liga de Expansion"
$ws.Cells.Item(217, 4).Value = "Mexico Liga de Expansion"
$ws.Cells.Item(217, 5).Value = 45384.83680555555
$ws.Cells.Item(217, 6).Value = "Tapatio"
$ws.Cells.Item(217, 7).Value = "Atlante"
$ws.Cells.Item(217, 8).Value = 0
$ws.Cells.Item(217, 9).Value = 0
$ws.Cells.Item(217, 10).Value = "D"
$ws.Cells.Item(217, 11).Value = 2.6
$ws.Cells.Item(217, 12).Value = 3.3
$ws.Cells.Item(217, 13).Value = 2.5
$ws.Cells.Item(217, 14).Value = 3.5
$ws.Cells.Item(217, 15).Value = 3.2
$ws.Cells.Item(217, 16).Value = 2.2
$ws.Cells.Item(217, 17).Value = 0.25
$ws.Cells.Item(217, 18).Value = 1.95
$ws.Cells.Item(217, 19).Value = 1.85
$ws.Cells.Item(217, 20).Value = 2.25
$ws.Cells.Item(217, 21).Value = 1.975
$ws.Cells.Item(217, 22).Value = 1.825
$ws.Cells.Item(217, 23).Value = -1
$ws.Cells.Item(217, 24).Value = 2.2
$ws.Cells.Item(217, 25).Value = -1
$ws.Cells.Item(217, 26).Value = 0.475
$ws.Cells.Item(217, 27).Value = -0.5
$ws.Cells.Item(217, 28).Value = -1
$ws.Cells.Item(217, 29).Value = 0.825

# Row 218
$ws.Cells.Item(218, 1).Value = 216
$ws.Cells.Item(218, 2).Value = 7641717
$ws.Cells.Item(218, 3).Value = "Mexico Liga de Expansion"
$ws.Cells.Item(218, 4).Value = "Mexico Liga de Expansion"
$ws.Cells.Item(218, 5).Value = 45384.92013888889
$ws.Cells.Item(218, 6).Value = "Cimarrones de Sonora FC"
$ws.Cells.Item(218, 7).Value = "Club Atletico La Paz"
$ws.Cells.Item(218, 8).Value = 2
$ws.Cells.Item(218, 9).Value = 1
$ws.Cells.Item(218, 10).Value = "H"
$ws.Cells.Item(218, 11).Value = 2.15
$ws.Cells.Item(218, 12).Value = 3.25
$ws.Cells.Item(218, 13).Value = 3.25
$ws.Cells.Item(218, 14).Value = 1.8
$ws.Cells.Item(218, 15).Value = 3.75
$ws.Cells.Item(218, 16).Value = 4.5
$ws.Cells.Item(218, 17).Value = -0.5
$ws.Cells.Item(218, 18).Value = 1.8
$ws.Cells.Item(218, 19).Value = 2
$ws.Cells.Item(218, 20).Value = 2.5
$ws.Cells.Item(218, 21).Value = 1.85
$ws.Cells.Item(218, 22).Value = 1.95
$ws.Cells.Item(218, 23).Value = 0.8
$ws.Cells.Item(218, 24).Value = -1
$ws.Cells.Item(218, 25).Value = -1
$ws.Cells.Item(218, 26).Value = 0.8
$ws.Cells.Item(218, 27).Value = -1
$ws.Cells.Item(218, 28).Value = 0.8500000000000001
$ws.Cells.Item(218, 29).Value = -1

# Row 219
$ws.Cells.Item(219, 1).Value = 217
$ws.Cells.Item(219, 2).Value = 7641718
$ws.Cells.Item(219, 3).Value = "Mexico Liga de Expansion"
$ws.Cells.Item(219, 4).Value = "Mexico Liga de Expansion"
$ws.Cells.Item(219, 5).Value = 45385.00347222222
$ws.Cells.Item(219, 6).Value = "Dorados"
$ws.Cells.Item(219, 7).Value = "Oaxaca"
$ws.Cells.Item(219, 8).Value = 1
$ws.Cells.Item(219, 9).Value = 3
$ws.Cells.Item(219, 10).Value = "A"
$ws.Cells.Item(219, 11).Value = 2.25
$ws.Cells.Item(219, 12).Value = 3.25
$ws.Cells.Item(219, 13).Value = 3
$ws.Cells.Item(219, 14).Value = 2.3
$ws.Cells.Item(219, 15).Value = 3.25
$ws.Cells.Item(219, 16).Value = 3.2
$ws.Cells.Item(219, 17).Value = -0.25
$ws.Cells.Item(219, 18).Value = 2
$ws.Cells.Item(219, 19).Value = 1.8
$ws.Cells.Item(219, 20).Value = 2.5
$ws.Cells.Item(219, 21).Value = 1.975
$ws.Cells.Item(219, 22).Value = 1.825
$ws.Cells.Item(219, 23).Value = -1
$ws.Cells.Item(219, 24).Value = -1
$ws.Cells.Item(219, 25).Value = 2.2
$ws.Cells.Item(219, 26).Value = -1
$ws.Cells.Item(219, 27).Value = 0.8
$ws.Cells.Item(219, 28).Value = 0.9750000000000001
$ws.Cells.Item(219, 29).Value = -1

# --- Append new rows 220-223 ---
# Copy cell formatting (bold/bordered id column, date format column, etc.)
# from an existing fully-populated row so new rows match existing look & feel.
$ws.Range("A216:AC216").Copy()
$ws.Range("A220:AC223").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Row 220 (new)
$ws.Cells.Item(220, 1).Value = 218
$ws.Cells.Item(220, 2).Value = 7641719
$ws.Cells.Item(220, 3).Value = "Mexico Liga de Expansion"
$ws.Cells.Item(220, 4).Value = "Mexico Liga de Expansion"
$ws.Cells.Item(220, 5).Value = 45385.92013888889
$ws.Cells.Item(220, 6).Value = "Mineros de Zacatecas"
$ws.Cells.Item(220, 7).Value = "Club Celaya"
$ws.Cells.Item(220, 8).Value = 2
$ws.Cells.Item(220, 9).Value = 2
$ws.Cells.Item(220, 10).Value = "D"
$ws.Cells.Item(220, 11).Value = 2.15
$ws.Cells.Item(220, 12).Value = 3.25
$ws.Cells.Item(220, 13).Value = 3.2
$ws.Cells.Item(220, 14).Value = 2.3
$ws.Cells.Item(220, 15).Value = 3.3
$ws.Cells.Item(220, 16).Value = 3.1
$ws.Cells.Item(220, 17).Value = -0.25
$ws.Cells.Item(220, 18).Value = 1.975
$ws.Cells.Item(220, 19).Value = 1.825
$ws.Cells.Item(220, 20).Value = 2.5
$ws.Cells.Item(220, 21).Value = 1.85
$ws.Cells.Item(220, 22).Value = 1.95
$ws.Cells.Item(220, 23).Value = -1
$ws.Cells.Item(220, 24).Value = 2.3
$ws.Cells.Item(220, 25).Value = -1
$ws.Cells.Item(220, 26).Value = -0.5
$ws.Cells.Item(220, 27).Value = 0.4125
$ws.Cells.Item(220, 28).Value = 0.8500000000000001
$ws.Cells.Item(220, 29).Value = -1

# Row 221 (new)
$ws.Cells.Item(221, 1).Value = 219
$ws.Cells.Item(221, 2).Value = 7641720
$ws.Cells.Item(221, 3).Value = "Mexico Liga de Expansion"
$ws.Cells.Item(221, 4).Value = "Mexico Liga de Expansion"
$ws.Cells.Item(221, 5).Value = 45386.00347222222
$ws.Cells.Item(221, 6).Value = "Venados FC"
$ws.Cells.Item(221, 7).Value = "Cancun FC"
$ws.Cells.Item(221, 8).Value = 2
$ws.Cells.Item(221, 9).Value = 0
$ws.Cells.Item(221, 10).Value = "H"
$ws.Cells.Item(221, 11).Value = 2.5
$ws.Cells.Item(221, 12).Value = 3.1
$ws.Cells.Item(221, 13).Value = 2.75
$ws.Cells.Item(221, 14).Value = 2.5
$ws.Cells.Item(221, 15).Value = 3.25
$ws.Cells.Item(221, 16).Value = 2.875
$ws.Cells.Item(221, 17).Value = 0
$ws.Cells.Item(221, 18).Value = 1.775
$ws.Cells.Item(221, 19).Value = 2.025
$ws.Cells.Item(221, 20).Value = 2.25
$ws.Cells.Item(221, 21).Value = 1.8
$ws.Cells.Item(221, 22).Value = 2
$ws.Cells.Item(221, 23).Value = 1.5
$ws.Cells.Item(221, 24).Value = -1
$ws.Cells.Item(221, 25).Value = -1
$ws.Cells.Item(221, 26).Value = 0.7749999999999999
$ws.Cells.Item(221, 27).Value = -1
$ws.Cells.Item(221, 28).Value = -0.5
$ws.Cells.Item(221, 29).Value = 0.5

# Row 222 (new)
$ws.Cells.Item(222, 1).Value = 220
$ws.Cells.Item(222, 2).Value = 7641721
$ws.Cells.Item(222, 3).Value = "Mexico Liga de Expansion"
$ws.Cells.Item(222, 4).Value = "Mexico Liga de Expansion"
$ws.Cells.Item(222, 5).Value = 45386.92013888889
$ws.Cells.Item(222, 6).Value = "Atletico Morelia"
$ws.Cells.Item(222, 7).Value = "Universidad Guadalajara"
$ws.Cells.Item(222, 8).Value = 0
$ws.Cells.Item(222, 9).Value = 2
$ws.Cells.Item(222, 10).Value = "A"
$ws.Cells.Item(222, 11).Value = 2.6
$ws.Cells.Item(222, 12).Value = 3.3
$ws.Cells.Item(222, 13).Value = 2.5
$ws.Cells.Item(222, 14).Value = 3.3
$ws.Cells.Item(222, 15).Value = 3.3
$ws.Cells.Item(222, 16).Value = 2.25
$ws.Cells.Item(222, 17).Value = 0.25
$ws.Cells.Item(222, 18).Value = 1.825
$ws.Cells.Item(222, 19).Value = 1.975
$ws.Cells.Item(222, 20).Value = 2.25
$ws.Cells.Item(222, 21).Value = 1.8
$ws.Cells.Item(222, 22).Value = 2
$ws.Cells.Item(222, 23).Value = -1
$ws.Cells.Item(222, 24).Value = -1
$ws.Cells.Item(222, 25).Value = 1.25
$ws.Cells.Item(222, 26).Value = -1
$ws.Cells.Item(222, 27).Value = 0.9750000000000001
$ws.Cells.Item(222, 28).Value = -0.5
$ws.Cells.Item(222, 29).Value = 0.5

# Row 223 (new)
$ws.Cells.Item(223, 1).Value = 221
$ws.Cells.Item(223, 2).Value = 7641722
$ws.Cells.Item(223, 3).Value = "Mexico Liga de Expansion"
$ws.Cells.Item(223, 4).Value = "Mexico Liga de Expansion"
$ws.Cells.Item(223, 5).Value = 45387.00347222222
$ws.Cells.Item(223, 6).Value = "Tepatitlan FC"
$ws.Cells.Item(223, 7).Value = "Tlaxcala FC"
$ws.Cells.Item(223, 8).Value = 1
$ws.Cells.Item(223, 9).Value = 1
$ws.Cells.Item(223, 10).Value = "D"
$ws.Cells.Item(223, 11).Value = 2.6
$ws.Cells.Item(223, 12).Value = 3.1
$ws.Cells.Item(223, 13).Value = 2.65
$ws.Cells.Item(223, 14).Value = 3
$ws.Cells.Item(223, 15).Value = 3.3
$ws.Cells.Item(223, 16).Value = 2.4
$ws.Cells.Item(223, 17).Value = 0
$ws.Cells.Item(223, 18).Value = 2
$ws.Cells.Item(223, 19).Value = 1.8
$ws.Cells.Item(223, 20).Value = 2.25
$ws.Cells.Item(223, 21).Value = 1.85
$ws.Cells.Item(223, 22).Value = 1.95
$ws.Cells.Item(223, 23).Value = -1
$ws.Cells.Item(223, 24).Value = 2.3
$ws.Cells.Item(223, 25).Value = -1
$ws.Cells.Item(223, 26).Value = 0
$ws.Cells.Item(223, 27).Value = -0
$ws.Cells.Item(223, 28).Value = -0.5
$ws.Cells.Item(223, 29).Value = 0.475

